# Divide la interfaz de reservas de restaurante: libera la mesa 2 (quita el
# usuario asignado y cambia su estado a "Libre") y agrega las mesas 5 y 10
# con su estado y capacidad.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mesa 2 (fila 2): quitar el usuario asignado y marcarla como Libre
$ws.Range("B2").ClearContents()
$ws.Range("C2").Value = "Libre"

# Mesa 5 (fila 4, nueva)
$ws.Range("A4").Value = 5
$ws.Range("C4").Value = "Libre"
$ws.Range("D4").Value = 4

# Mesa 10 (fila 5, nueva)
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "'"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = "Libre"
$ws.Range("D5").Value = "'2"
$ws.Range("D5").ClearFormats()
